$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1404.8182
$ws.Range("J17").Value = 1404.8182
$ws.Range("L17").Value = 4214.4546
$ws.Range("N17").Value = -4550.4546
$ws.Range("H51").Value = 4940.1816
$ws.Range("I51").Value = 5123.8335
$ws.Range("J51").Value = 4719.8
$ws.Range("K51").Value = 5123.8335
$ws.Range("L51").Value = 4719.8
$ws.Range("M51").Value = -4639.8335
$ws.Range("N51").Value = -5687.8
$ws.Range("H55").Value = 832.6667
$ws.Range("I55").Value = 498
$ws.Range("J55").Value = 1000
$ws.Range("K55").Value = 498
$ws.Range("L55").Value = 1000
$ws.Range("M55").Value = -284
$ws.Range("N55").Value = -1428
$ws.Range("H92").Value = 20433.4
$ws.Range("I92").Value = 25486.4
$ws.Range("J92").Value = 221.4
$ws.Range("K92").Value = 25486.4
$ws.Range("L92").Value = 221.4
$ws.Range("M92").Value = -24238.4
$ws.Range("N92").Value = -2717.4
$ws.Range("H98").Value = 1791.0869
$ws.Range("I98").Value = 1922.6111
$ws.Range("J98").Value = 1317.6
$ws.Range("K98").Value = 1922.6111
$ws.Range("L98").Value = 1317.6
$ws.Range("M98").Value = -424.6111000000001
$ws.Range("N98").Value = -4313.6
$ws.Range("H116").Value = 7739.304
$ws.Range("I116").Value = 7444.3335
$ws.Range("J116").Value = 8801.200000000001
$ws.Range("K116").Value = 7444.3335
$ws.Range("L116").Value = 8801.200000000001
$ws.Range("M116").Value = -4002.3335
$ws.Range("N116").Value = -15685.2
$ws.Range("H122").Value = 1791.0869
$ws.Range("I122").Value = 1922.6111
$ws.Range("J122").Value = 1317.6
$ws.Range("K122").Value = 5767.8333
$ws.Range("L122").Value = 3952.8
$ws.Range("M122").Value = -3317.8333
$ws.Range("N122").Value = -8852.799999999999
$ws.Range("H125").Value = 8838.429
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
$ws.Range("H132").Value = 1823.2808
$ws.Range("I132").Value = 1575.9811
$ws.Range("K132").Value = 4727.9433
$ws.Range("M132").Value = -2197.9433

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1634.84
$ws.Range("I4").Value = 1327.75
$ws.Range("J4").Value = 2180.7778
$ws.Range("K4").Value = 1327.75
$ws.Range("L4").Value = 2180.7778
$ws.Range("M4").Value = -1211.75
$ws.Range("N4").Value = -2412.7778
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("H61").Value = 32194.938
$ws.Range("I61").Value = 2413.111
$ws.Range("K61").Value = 2413.111
$ws.Range("M61").Value = -2201.111
$ws.Range("H63").Value = 2802.2
$ws.Range("I63").Value = 2502.75
$ws.Range("J63").Value = 4000
$ws.Range("K63").Value = 2502.75
$ws.Range("L63").Value = 4000
$ws.Range("M63").Value = -1816.75
$ws.Range("N63").Value = -5372
$ws.Range("H66").Value = 2802.2
$ws.Range("I66").Value = 2502.75
$ws.Range("J66").Value = 4000
$ws.Range("K66").Value = 12513.75
$ws.Range("L66").Value = 20000
$ws.Range("M66").Value = -9081.75
$ws.Range("N66").Value = -26864
$ws.Range("H74").Value = 258954.14
$ws.Range("I74").Value = 226191.5
$ws.Range("K74").Value = 226191.5
$ws.Range("M74").Value = -225317.5
$ws.Range("H77").Value = 258954.14
$ws.Range("I77").Value = 226191.5
$ws.Range("K77").Value = 1130957.5
$ws.Range("M77").Value = -1126589.5
$ws.Range("H97").Value = 630.7646999999999
$ws.Range("I97").Value = 630.7646999999999
$ws.Range("K97").Value = 630.7646999999999
$ws.Range("M97").Value = -134.7646999999999
$ws.Range("H132").Value = 7223.38
$ws.Range("I132").Value = 7531.087
$ws.Range("K132").Value = 22593.261
$ws.Range("M132").Value = -20063.261
$ws.Range("H136").Value = 32194.938
$ws.Range("I136").Value = 2413.111
$ws.Range("K136").Value = 7239.333
$ws.Range("M136").Value = -4689.333
$ws.Range("M6").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 29999
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1831.0209
$ws.Range("I31").Value = 1230.1177
$ws.Range("J31").Value = 3290.3572
$ws.Range("K31").Value = 1230.1177
$ws.Range("L31").Value = 3290.3572
$ws.Range("M31").Value = -935.1177
$ws.Range("N31").Value = -3880.3572
$ws.Range("H34").Value = 1831.0209
$ws.Range("I34").Value = 1230.1177
$ws.Range("J34").Value = 3290.3572
$ws.Range("K34").Value = 1230.1177
$ws.Range("L34").Value = 3290.3572
$ws.Range("M34").Value = -1028.1177
$ws.Range("N34").Value = -3694.3572
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H122").Value = 2338.2
$ws.Range("I122").Value = 2199
$ws.Range("K122").Value = 6597
$ws.Range("M122").Value = -4147
$ws.Range("H132").Value = 3824.1064
$ws.Range("I132").Value = 1864.6316
$ws.Range("K132").Value = 5593.8948
$ws.Range("M132").Value = -3063.8948
$ws.Range("H134").Value = 3138.3125
$ws.Range("I134").Value = 2631.2856
$ws.Range("K134").Value = 7893.8568
$ws.Range("M134").Value = -5358.8568
$ws.Range("M76").ClearContents()
$ws.Range("M79").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 33598.832
$ws.Range("I9").Value = 33598.832
$ws.Range("K9").Value = 100796.496
$ws.Range("M9").Value = -100572.496
$ws.Range("H60").Value = 1798.3077
$ws.Range("I60").Value = 1269
$ws.Range("K60").Value = 3807
$ws.Range("M60").Value = -3556
$ws.Range("H137").Value = 4388.19
$ws.Range("I137").Value = 1087.75
$ws.Range("J137").Value = 4838.25
$ws.Range("K137").Value = 3263.25
$ws.Range("L137").Value = 14514.75
$ws.Range("M137").Value = 1836.75
$ws.Range("N137").Value = -24714.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 2505555
$ws.Range("I12").Value = 4999999
$ws.Range("J12").Value = 11111
$ws.Range("K12").Value = 4999999
$ws.Range("L12").Value = 11111
$ws.Range("M12").Value = -4999859
$ws.Range("N12").Value = -11391
$ws.Range("H18").Value = 8699.75
$ws.Range("I18").Value = 4931
$ws.Range("K18").Value = 4931
$ws.Range("M18").Value = -4638
$ws.Range("H102").Value = 33854.676
$ws.Range("I102").Value = 40834.777
$ws.Range("K102").Value = 40834.777
$ws.Range("M102").Value = -39212.777
$ws.Range("H109").Value = 31749
$ws.Range("J109").Value = 31749
$ws.Range("L109").Value = 31749
$ws.Range("N109").Value = -33829
$ws.Range("H113").Value = 4903.3335
$ws.Range("I113").Value = 5701.25
$ws.Range("J113").Value = 2350
$ws.Range("K113").Value = 5701.25
$ws.Range("L113").Value = 2350
$ws.Range("M113").Value = -3531.25
$ws.Range("N113").Value = -6690
$ws.Range("H132").Value = 3558.5
$ws.Range("I132").Value = 3777.4443
$ws.Range("K132").Value = 11332.3329
$ws.Range("M132").Value = -8802.332900000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8734.799999999999
$ws.Range("I40").Value = 7418.625
$ws.Range("K40").Value = 7418.625
$ws.Range("M40").Value = -7282.625
$ws.Range("H125").Value = 79599
$ws.Range("J125").Value = 79599
$ws.Range("L125").Value = 79599
$ws.Range("N125").Value = -89439
$ws.Range("H132").Value = 4443.7085
$ws.Range("I132").Value = 3024.7222
$ws.Range("K132").Value = 9074.1666
$ws.Range("M132").Value = -6544.1666
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3692.8572
$ws.Range("I132").Value = 1707.5714
$ws.Range("J132").Value = 11634
$ws.Range("K132").Value = 5122.7142
$ws.Range("L132").Value = 34902
$ws.Range("M132").Value = -2592.7142
$ws.Range("N132").Value = -39962
$ws.Range("H133").Value = 69398
$ws.Range("J133").Value = 69398
$ws.Range("L133").Value = 69398
$ws.Range("N133").Value = -79518
